$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (column F) values for the affected rows.
$updates = @{
    2  = -4
    4  = -8
    7  = -12
    10 = -4
    11 = -6
    12 = -7
    13 = -11
    14 = -12
    15 = -11
    17 = -5
    20 = 2
    22 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
